# Swap the species-record data between row 3 and row 4 on the active sheet.
# Only the columns that actually differ between the two rows are touched:
# A, B, D, E, F, G, H, Q, R, Z, AB

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")

foreach ($col in $columns) {
    $addr3 = "{0}3" -f $col
    $addr4 = "{0}4" -f $col

    $cell3 = $ws.Range($addr3)
    $cell4 = $ws.Range($addr4)

    $value3 = $cell3.Value2
    $value4 = $cell4.Value2

    $cell3.Value2 = $value4
    $cell4.Value2 = $value3
}
